$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.019613333333333
$ws.Range("H2").Value = 12.05884
$ws.Range("I2").Value = 0.4377217086785624
$ws.Range("J2").Value = 0.4377217086785624
$ws.Range("M2").Value = 4.019613333333333
$ws.Range("N2").Value = 12.05884
$ws.Range("O2").Value = 0.4377217086785624
$ws.Range("P2").Value = 0.4377217086785624
$ws.Range("Q2").Value = 16.15729134951111
$ws.Range("R2").Value = 145.4156221456
$ws.Range("S2").Value = 0.1916002942484802
$ws.Range("T2").Value = 0.1916002942484802
$ws.Range("G3").Value = 4.019613333333333
$ws.Range("H3").Value = 12.05884
$ws.Range("I3").Value = 0.4377217086785624
$ws.Range("J3").Value = 0.4377217086785624
$ws.Range("O3").Value = 0.02575118419467902
$ws.Range("P3").Value = 0.02575118419467902
$ws.Range("Q3").Value = 0.9505340433866666
$ws.Range("R3").Value = 8.55480639048
$ws.Range("S3").Value = 0.01127185234619129
$ws.Range("T3").Value = 0.01127185234619129
$ws.Range("G4").Value = 4.019613333333333
$ws.Range("H4").Value = 12.05884
$ws.Range("I4").Value = 0.4377217086785624
$ws.Range("J4").Value = 0.4377217086785624
$ws.Range("M4").Value = 4.926946666666667
$ws.Range("N4").Value = 14.78084
$ws.Range("O4").Value = 0.5365271071267587
$ws.Range("P4").Value = 0.5365271071267587
$ws.Range("Q4").Value = 19.80442051395556
$ws.Range("R4").Value = 178.2397846256
$ws.Range("S4").Value = 0.2348495620838909
$ws.Range("T4").Value = 0.2348495620838909
$ws.Range("I5").Value = 0.02575118419467902
$ws.Range("J5").Value = 0.02575118419467902
$ws.Range("M5").Value = 4.019613333333333
$ws.Range("N5").Value = 12.05884
$ws.Range("O5").Value = 0.4377217086785624
$ws.Range("P5").Value = 0.4377217086785624
$ws.Range("Q5").Value = 0.9505340433866666
$ws.Range("R5").Value = 8.55480639048
$ws.Range("S5").Value = 0.01127185234619129
$ws.Range("T5").Value = 0.01127185234619129
$ws.Range("I6").Value = 0.02575118419467902
$ws.Range("J6").Value = 0.02575118419467902
$ws.Range("O6").Value = 0.02575118419467902
$ws.Range("P6").Value = 0.02575118419467902
$ws.Range("S6").Value = 0.0006631234874282864
$ws.Range("T6").Value = 0.0006631234874282864
$ws.Range("I7").Value = 0.02575118419467902
$ws.Range("J7").Value = 0.02575118419467902
$ws.Range("M7").Value = 4.926946666666667
$ws.Range("N7").Value = 14.78084
$ws.Range("O7").Value = 0.5365271071267587
$ws.Range("P7").Value = 0.5365271071267587
$ws.Range("Q7").Value = 1.165094786053333
$ws.Range("R7").Value = 10.48585307448
$ws.Range("S7").Value = 0.01381620836105944
$ws.Range("T7").Value = 0.01381620836105944
$ws.Range("G8").Value = 4.926946666666667
$ws.Range("H8").Value = 14.78084
$ws.Range("I8").Value = 0.5365271071267587
$ws.Range("J8").Value = 0.5365271071267587
$ws.Range("M8").Value = 4.019613333333333
$ws.Range("N8").Value = 12.05884
$ws.Range("O8").Value = 0.4377217086785624
$ws.Range("P8").Value = 0.4377217086785624
$ws.Range("Q8").Value = 19.80442051395556
$ws.Range("R8").Value = 178.2397846256
$ws.Range("S8").Value = 0.2348495620838909
$ws.Range("T8").Value = 0.2348495620838909
$ws.Range("G9").Value = 4.926946666666667
$ws.Range("H9").Value = 14.78084
$ws.Range("I9").Value = 0.5365271071267587
$ws.Range("J9").Value = 0.5365271071267587
$ws.Range("O9").Value = 0.02575118419467902
$ws.Range("P9").Value = 0.02575118419467902
$ws.Range("Q9").Value = 1.165094786053333
$ws.Range("R9").Value = 10.48585307448
$ws.Range("S9").Value = 0.01381620836105944
$ws.Range("T9").Value = 0.01381620836105944
$ws.Range("G10").Value = 4.926946666666667
$ws.Range("H10").Value = 14.78084
$ws.Range("I10").Value = 0.5365271071267587
$ws.Range("J10").Value = 0.5365271071267587
$ws.Range("M10").Value = 4.926946666666667
$ws.Range("N10").Value = 14.78084
$ws.Range("O10").Value = 0.5365271071267587
$ws.Range("P10").Value = 0.5365271071267587
$ws.Range("Q10").Value = 24.27480345617778
$ws.Range("R10").Value = 218.4732311056
$ws.Range("S10").Value = 0.2878613366818084
$ws.Range("T10").Value = 0.2878613366818084
